$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: politeness_score becomes a real number (was stored as text "3")
$ws.Range("B6").Value = 3

# Row 7: new annotation row appended below the existing data
$ws.Range("A7").Value = "Sunsi Wu"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "3"
$ws.Range("C7").Value = "did"
$ws.Range("D7").Value = "FBK"
$ws.Range("E7").Value = "RES"
$ws.Range("F7").Value = "dc9804e9-fe90-49ab-88bb-ac97478c1b97"
$ws.Range("G7").Value = "i87JIQTAnB8AQ_annotated.xlsx"
$ws.Range("H7").Value = "As you suggested, I did run comparison tests and I will present the results here."
